$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64, shifting existing rows 64-68 down to 65-69
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new data entry
$ws.Range("A64").Value = 11
$ws.Range("B64").Value = "Vega Monumental Concepción"
$ws.Range("C64").Value = "Bíobío"
$ws.Range("D64").Value = 44946
$ws.Range("D64").NumberFormat = $ws.Range("D65").NumberFormat
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = 100112031
$ws.Range("G64").Value = "Poroto verde"
$ws.Range("H64").Value = "Magnum"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 150
$ws.Range("K64").Value = 16000
$ws.Range("L64").Value = 17000
$ws.Range("M64").Value = 16533
$ws.Range("N64").Value = "$/saco 25 kilos"
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 661
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"
